# Generate Report for Handback
# Marks the zh-cn / de-de localization rows as handed back: updates the
# "Status" text (Overview summary + per-language sheets), records the
# source file as the "Latest Target File" (with a hyperlink back to the
# source .md, same as column A), fills in the "Latest Handback File" with
# the generated .xlf name, stamps the "Latest Handback DateTime", and
# widens a few columns so the new/longer values aren't clipped.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$sourceMd   = "8b10ff1d-7258-479b-9cee-88ff74bfa152.md"
$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/8b10ff1d-7258-479b-9cee-88ff74bfa152.md"

# ---- Overview sheet: zh-cn / de-de status summary columns ----
$ws_overview.Range("E2").Value = $statusText
$ws_overview.Range("F2").Value = $statusText
$ws_overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$ws_overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---- zh-cn sheet ----
$ws_zhcn.Range("C2").Value = $statusText
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $sourceMdUrl, $null, $sourceMd, $sourceMd)
$ws_zhcn.Range("J2").Value = "8b10ff1d-7258-479b-9cee-88ff74bfa152.7929c32eb86e154314d1e8ead1699068515dae77.zh-cn.xlf"
$ws_zhcn.Range("K2").Value = "2016-08-16 18:53:31"
$ws_zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws_zhcn.Columns.Item(9).ColumnWidth = 39.1666666666667
$ws_zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---- de-de sheet ----
$ws_dede.Range("C2").Value = $statusText
$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $sourceMdUrl, $null, $sourceMd, $sourceMd)
$ws_dede.Range("J2").Value = "8b10ff1d-7258-479b-9cee-88ff74bfa152.7929c32eb86e154314d1e8ead1699068515dae77.de-de.xlf"
$ws_dede.Range("K2").Value = "2016-08-16 18:53:38"
$ws_dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws_dede.Columns.Item(9).ColumnWidth = 39.1666666666667
$ws_dede.Columns.Item(10).ColumnWidth = 39.1666666666667
